$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Prepend 8 new paragraphs (team/title/milestone header block)
#    before the documents current first paragraph.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(1).Range
$anchor.Collapse(1)
for ($i = 0; $i -lt 8; $i++) {
    $anchor.InsertParagraphBefore()
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace each freshly-inserted (empty) paragraph with its final OOXML
# so the run formatting/paragraph formatting matches exactly.
$p1 = $d.Paragraphs.Item(1).Range
$p1.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Nhóm 23 </w:t></w:r></w:p>' + $xmlFooter)

$p2 = $d.Paragraphs.Item(2).Range
$p2.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Đề tài: Phát triển service xử lý và lưu trữ video</w:t></w:r></w:p>' + $xmlFooter)

$p3 = $d.Paragraphs.Item(3).Range
$p3.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>' + $xmlFooter)

$p4 = $d.Paragraphs.Item(4).Range
$p4.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Họ và tên: Nguyễn Minh Sơn</w:t></w:r></w:p>' + $xmlFooter)

$p5 = $d.Paragraphs.Item(5).Range
$p5.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>MSSV: 20110713</w:t></w:r></w:p>' + $xmlFooter)

$p6 = $d.Paragraphs.Item(6).Range
$p6.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>' + $xmlFooter)

$p7 = $d.Paragraphs.Item(7).Range
$p7.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Milestone 2</w:t></w:r></w:p>' + $xmlFooter)

$p8 = $d.Paragraphs.Item(8).Range
$p8.InsertXML($xmlHeader + '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>' + $xmlFooter)

# ------------------------------------------------------------------
# 2) Add the "Normal (Web)" paragraph style definition (styles.xml).
#    It mirrors what Word auto-adds when web/HTML content is pasted
#    in; it is not referenced by any paragraph in the document.
# ------------------------------------------------------------------
$normalWeb = $d.Styles.Add("Normal (Web)", 1)
$normalWeb.BaseStyle = $d.Styles.Item("Normal")
$normalWeb.Priority = 99
$normalWeb.UnhideWhenUsed = $true

$normalWeb.Font.Name = "Times New Roman"
$normalWeb.Font.NameFarEast = "Times New Roman"
$normalWeb.Font.NameBi = "Times New Roman"
$normalWeb.Font.Size = 12
$normalWeb.Font.SizeBi = 12

$normalWeb.ParagraphFormat.SpaceBefore = 5
$normalWeb.ParagraphFormat.SpaceBeforeAuto = $true
$normalWeb.ParagraphFormat.SpaceAfter = 5
$normalWeb.ParagraphFormat.SpaceAfterAuto = $true
$normalWeb.ParagraphFormat.LineSpacingRule = 0

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Paragraph 9 (former first paragraph) text:" $d.Paragraphs.Item(9).Range.Text
